# Generate Report for Handoff
# Update the GUID-based file names and handoff/handback timestamps
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "4c997830-cf2f-4dd3-831a-bf512d97ef0f"
$newGuid = "6cbd9040-2d4e-4c40-8a58-d6278cda3a29"

$oldHash = "6b45f7c4424769db43d31dcb68cbb18564a8eda6"
$newHash = "662766103e86beab5d4692c118b97387b2b6c8f0"

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-03-24 02:52:41"

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-24 02:52:33"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-24 02:52:41"
